$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header label above the statistical error block
$ws.Range("I1").Value = "Statistical errors"

# Column headers for the statistical error values (bold, like A5:E5)
$ws.Range("I3").Value = "60h"
$ws.Range("J3").Value = "HighKick"
$ws.Range("K3").Value = "9d"
$ws.Range("L3").Value = "Endgame"
$ws.Range("I3:L3").Font.Bold = $true

# Statistical error values (row 5 otherwise carries an inherited bold row
# format, so reset the style on these cells back to the workbook default)
$ws.Range("I5").Style = $ws.Range("H1").Style
$ws.Range("J5").Style = $ws.Range("H1").Style
$ws.Range("K5").Style = $ws.Range("H1").Style
$ws.Range("L5").Style = $ws.Range("H1").Style
$ws.Range("I5").Value = 1358.1
$ws.Range("J5").Value = 1411.2
$ws.Range("K5").Value = 903.3
$ws.Range("L5").Value = 639.3

# Label for total errors block
$ws.Range("I8").Value = "Total errors"

# Column headers for the total error values
$ws.Range("I10").Value = "60h"
$ws.Range("J10").Value = "HighKick"
$ws.Range("K10").Value = "9d"
$ws.Range("L10").Value = "Endgame"
$ws.Range("I10:L10").Font.Bold = $true

# Total error formulas (quadrature sum of statistical error and total quadrature sum)
$ws.Range("I12:L12").Formula = "=SQRT(SUMSQ(I5,B33))"

# Column widths
$ws.Columns.Item(9).ColumnWidth = 16.1640625
$ws.Columns.Item(12).ColumnWidth = 11.83203125

$ws.Range("I17").Select()

$wb.Save()
